# Backlog.xlsx — "FW requirements" sheet: record the assignee ("Wyatt Wang")
# for the first three backlog items, and leave the selection on D8 (the last
# cell touched), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FW requirements")
$ws.Activate() | Out-Null

$ws.Range("F2").Value = "Wyatt Wang"
$ws.Range("F3").Value = "Wyatt Wang"
$ws.Range("F4").Value = "Wyatt Wang"

$ws.Range("D8").Select() | Out-Null
